$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 288; this shifts the existing rows 288..362
# down to 289..363 (carrying their values/formatting with them, matching
# the target diff where old row 288 becomes row 289, ..., old row 362
# becomes row 363).
$ws.Rows("288").Insert()

# Populate the newly inserted row 288 with its data. Columns A, B, C, E, F,
# G, H, I, J and R repeat the same values as the (now shifted) row below
# (old row 288, now row 289), while D, K, L, M, N, O, P, Q hold new values.
$ws.Range("A288").Value = 5
$ws.Range("B288").Value = "Macroferia Regional de Talca"
$ws.Range("C288").Value = "Maule"
$ws.Range("D288").Value = 44722
$ws.Range("E288").Value = 7
$ws.Range("F288").Value = 100112032
$ws.Range("G288").Value = "Zapallo italiano"
$ws.Range("H288").Value = "Sin especificar"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 300
$ws.Range("K288").Value = 8000
$ws.Range("L288").Value = 8000
$ws.Range("M288").Value = 8000
$ws.Range("N288").Value = "$/caja 50 unidades"
$ws.Range("O288").Value = "Región de Arica y Parinacota"
$ws.Range("P288").Value = 160
$ws.Range("Q288").Value = 50
$ws.Range("R288").Value = "Hortaliza"
